# Select grouping level and groups for time series.
#
# Adds the new VR5x / VR4 hydrography station rows (40-45) to the KML
# sheet, then leaves the selection where the user ended up after typing
# them in (E46).
#
# The writes below are ordered to match the exact sequence the data was
# originally typed/pasted into the sheet (column A down, then G/H filled
# down from the row above, then C back-filled, then A/B for the next
# block of rows, etc.) so that new shared-string entries are interned in
# the same order as the source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KML")

$ws.Cells.Item(40, 1).Value = "VR55"
$ws.Cells.Item(40, 2).Value = "Spilderbukta"

$ws.Cells.Item(40, 4).Value = 69.9664
$ws.Cells.Item(40, 5).Value = 21.6887
$ws.Cells.Item(40, 6).Value = 0

$ws.Cells.Item(40, 7).Value = "NA"
$ws.Cells.Item(40, 8).Value = "NA"
$ws.Cells.Item(41, 7).Value = "NA"
$ws.Cells.Item(41, 8).Value = "NA"
$ws.Cells.Item(42, 7).Value = "NA"
$ws.Cells.Item(42, 8).Value = "NA"
$ws.Cells.Item(43, 7).Value = "NA"
$ws.Cells.Item(43, 8).Value = "NA"
$ws.Cells.Item(44, 7).Value = "NA"
$ws.Cells.Item(44, 8).Value = "NA"

$ws.Cells.Item(40, 3).Value = "Kvænangen"
$ws.Cells.Item(45, 3).Value = "Kvænangen"

$ws.Cells.Item(41, 1).Value = "VR56"
$ws.Cells.Item(42, 1).Value = "VR57"
$ws.Cells.Item(42, 2).Value = "Storbukta"
$ws.Cells.Item(43, 1).Value = "VR58"
$ws.Cells.Item(43, 2).Value = "Ullsfjorden"
$ws.Cells.Item(44, 1).Value = "VR59"

$ws.Cells.Item(41, 2).Value = "Reisafjorden ytre"
$ws.Cells.Item(41, 3).Value = "Reisafjorden ytre"
$ws.Cells.Item(42, 3).Value = "Reisafjorden indre"
$ws.Cells.Item(43, 3).Value = "Ullsfjorden/Fugløyfjorden"
$ws.Cells.Item(44, 2).Value = "Sørfjorden ytre"
$ws.Cells.Item(44, 3).Value = "Sørfjorden ytre"

$ws.Cells.Item(41, 4).Value = 69.9068
$ws.Cells.Item(41, 5).Value = 21.0927
$ws.Cells.Item(41, 6).Value = 0

$ws.Cells.Item(42, 4).Value = 69.8515
$ws.Cells.Item(42, 5).Value = 21.1968
$ws.Cells.Item(42, 6).Value = 0

$ws.Cells.Item(43, 4).Value = 69.7544
$ws.Cells.Item(43, 5).Value = 19.7701
$ws.Cells.Item(43, 6).Value = 0

$ws.Cells.Item(44, 4).Value = 69.5711
$ws.Cells.Item(44, 5).Value = 19.7185
$ws.Cells.Item(44, 6).Value = 0

$ws.Cells.Item(45, 1).Value = "VR4"
$ws.Cells.Item(45, 2).Value = "Kvænangen ytre "

$ws.Cells.Item(45, 4).Value = 70.1161
$ws.Cells.Item(45, 5).Value = 21.0725

# Reflect the final selection/view after entering the new rows.
$ws.Range("E46").Select()
